$d = $word.ActiveDocument

# --- Bullet 1: "Utilize a Teensy 3.2-based interface to deliver precisely timed
#     digital pulses to initiate frame capture using a sCMOS camera"
#     -> "Developed a low-cost, open source, Teensy 3.2-based interface for
#     systems neuroscience experiments that can flexibly integrate diverse instruments"
$r = $d.Content
$r.Find.Execute(
    "Utilize a Teensy 3.2-based interface to deliver precisely timed digital pulses to initiate frame capture using a sCMOS camera",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Developed a low-cost, open source, Teensy 3.2-based interface for systems neuroscience experiments that can flexibly integrate diverse instruments",
    2) | Out-Null

# --- Bullet 2: "Demonstrate temporally precise behavioral data acquisition using a
#     Teensy 3.2 interface combined with concurrent sCMOS camera control"
#     -> "Demonstrate the temporal precision of the Teensy 3.2 interface in two
#     experimental conditions"
#     Done in two Find/Replace calls so the embedded _GoBack bookmark
#     (sitting between "Demonstrate" and "combined with") is left untouched.
$r = $d.Content
$r.Find.Execute(
    "Demonstrate temporally precise behavioral data acquisition using a Teensy 3.2 interface combined with ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Demonstrate",
    2) | Out-Null

$r = $d.Content
$r.Find.Execute(
    "concurrent sCMOS camera control",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " the temporal precision of the Teensy 3.2 interface in two experimental conditions",
    2) | Out-Null

# --- Bullet 3: "Demonstrate experimental control of combined analog (sound waveforms)
#     and digital pulses delivered simultaneously with camera control"
#     -> "Demonstrate the utility of the Teensy 3.2 interface in generating analog
#     signals (sound waveforms) and digital pulses simultaneously"
$r = $d.Content
$r.Find.Execute(
    "experimental control of combined analog (sound waveforms) and digital pulses delivered simultaneously with camera control",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the utility of the Teensy 3.2 interface in generating analog signals (sound waveforms) and digital pulses simultaneously",
    2) | Out-Null

# --- Remove the trailing empty paragraph at the end of the document.
$lastIdx = $d.Paragraphs.Count
$secondLastIdx = $lastIdx - 1
$pLast = $d.Paragraphs.Item($lastIdx)
$pPrev = $d.Paragraphs.Item($secondLastIdx)
$delRange = $d.Range($pPrev.Range.End - 1, $pLast.Range.End)
$delRange.Delete()
